$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.851.83'
$ws.Range('E2').Value = '  -1.11%  '

# Row 3
$ws.Range('E3').Value = '  -0.78%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7899'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.86%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.66'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.25%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3151'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.85%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.35'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.81%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07220'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.07%  '

# Row 11
$ws.Range('E11').Value = '  -0.15%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7659'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.16%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.515'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.25%  '

# Row 14
$ws.Range('D14').Value = '1.926.16'
$ws.Range('E14').Value = '  +0.81%  '

# Row 15
$ws.Range('E15').Value = '  -0.42%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.156'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.14%  '

# Row 17
$ws.Range('D17').Value = '29.875.85'
$ws.Range('E17').Value = '  -1.05%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.94'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.97%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.18'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.65%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007792'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.03%  '

# Row 21
$ws.Range('D21').Value = '2.173.22'
$ws.Range('E21').Value = '  +0.59%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.148'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +15.60%  '

# Row 24
$ws.Range('E24').Value = '  -0.12%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1650'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.53%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.414'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.87%  '

# Row 27
$ws.Range('E27').Value = '  -1.71%  '

# Row 28
$ws.Range('E28').Value = '  -1.56%  '

# Row 29
$ws.Range('E29').Value = '  -2.48%  '

# Row 30
$ws.Range('E30').Value = '  +2.46%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.549'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.63%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.506'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.58%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.114'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.58%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05569'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -8.62%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.272'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.31%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7439'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.35%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.40%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.614'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.65%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01926'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.56%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.775'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.61%  '

# Row 41
$ws.Range('D41').Value = '1.146.12'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.07'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.00%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4430'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.92%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.881'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.28%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8512'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.64%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.35'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.85%  '

# Row 47
$ws.Range('E47').Value = '  -0.03%  '

# Row 48
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.882'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.46%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.03'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.87%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.470'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.61%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.018'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +10.01%  '
